$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.942.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.318.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.91%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "187.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "554.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.310.29"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.578"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.97%  "
$ws.Range("E10").Value = "  -3.95%  "
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.90"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000265"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.848.66"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.45"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "579.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.928.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.338.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.62%  "
$ws.Range("E19").Value = "  +0.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.892"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "18.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.44%  "
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("E27").Value = "  +0.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.31"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "570.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.72"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "10.84"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.28%  "
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.734.22"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.95%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "55.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.20%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "34.07"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.04%  "
$ws.Range("E41").Value = "  -2.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₃0689"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.64%  "
$ws.Range("E44").Value = "  -7.14%  "
$ws.Range("E45").Value = "  +2.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.334"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0408"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.52%  "
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.16%  "
